# repull data, push all data, mean calculation
# Updates the "dSF" column (column F) values for several rows to reflect
# the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -12
    4  = 2
    6  = -3
    7  = 7
    9  = 1
    10 = -4
    11 = -1
    14 = -3
    16 = -4
    18 = 3
    19 = -7
    22 = -9
    23 = -5
    24 = 3
    25 = 1
    28 = -1
    29 = -1
    30 = 1
    32 = 1
    33 = 1
    34 = -3
    37 = 2
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
